$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 126.333336
$ws.Range("I4").Value = 111.2
$ws.Range("K4").Value = 111.2
$ws.Range("M4").Value = 2.799999999999997
# Row 15
$ws.Range("H15").Value = 116.35
$ws.Range("I15").Value = 116.35
$ws.Range("K15").Value = 349.05
$ws.Range("M15").Value = -180.05
# Row 17
$ws.Range("H17").Value = 3355492.2
$ws.Range("J17").Value = 3471133.5
$ws.Range("L17").Value = 10413400.5
$ws.Range("N17").Value = -10413736.5
# Row 54
$ws.Range("H54").Value = 9970
$ws.Range("I54").Value = 9925
$ws.Range("J54").Value = 10000
$ws.Range("K54").Value = 9925
$ws.Range("L54").Value = 10000
$ws.Range("M54").Value = -9439
$ws.Range("N54").Value = -10972
# Row 116
$ws.Range("H116").Value = 15684881
$ws.Range("I116").Value = 47044428
$ws.Range("J116").Value = 5108.5
$ws.Range("K116").Value = 47044428
$ws.Range("L116").Value = 5108.5
$ws.Range("M116").Value = -47040986
$ws.Range("N116").Value = -11992.5
# Row 118
$ws.Range("H118").Value = 460.63635
$ws.Range("I118").Value = 386.8
$ws.Range("J118").Value = 1199
$ws.Range("K118").Value = 1160.4
$ws.Range("L118").Value = 3597
$ws.Range("M118").Value = 496.5999999999999
$ws.Range("N118").Value = -6911
# Row 127
$ws.Range("H127").Value = 1351.375
$ws.Range("J127").Value = 2103.5
$ws.Range("L127").Value = 6310.5
$ws.Range("N127").Value = -16230.5
# Row 129
$ws.Range("H129").Value = 271198.62
$ws.Range("J129").Value = 295104.47
$ws.Range("L129").Value = 885313.4099999999
$ws.Range("N129").Value = -895313.4099999999
# Row 137
$ws.Range("H137").Value = 116154.31
$ws.Range("I137").Value = 155343.03
$ws.Range("J137").Value = 2942.4443
$ws.Range("K137").Value = 466029.09
$ws.Range("L137").Value = 8827.332900000001
$ws.Range("M137").Value = -463479.09
$ws.Range("N137").Value = -13927.3329
# Row 138
$ws.Range("H138").Value = 4077.5386
$ws.Range("I138").Value = 3029.5454
$ws.Range("J138").Value = 4846.067
$ws.Range("K138").Value = 9088.636200000001
$ws.Range("L138").Value = 14538.201
$ws.Range("M138").Value = -3948.636200000001
$ws.Range("N138").Value = -24818.201
# Row 141
$ws.Range("H141").Value = 1791.027
$ws.Range("I141").Value = 1653.5143
$ws.Range("K141").Value = 4960.5429
$ws.Range("M141").Value = 219.4570999999996

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1130.2084
$ws.Range("I2").Value = 1095.8235
$ws.Range("J2").Value = 1213.7142
$ws.Range("K2").Value = 1095.8235
$ws.Range("L2").Value = 1213.7142
$ws.Range("M2").Value = -982.8235
$ws.Range("N2").Value = -1439.7142
# Row 32
$ws.Range("H32").Value = 10289.102
$ws.Range("I32").Value = 7719.385
$ws.Range("J32").Value = 22219.928
$ws.Range("K32").Value = 7719.385
$ws.Range("L32").Value = 22219.928
$ws.Range("M32").Value = -7432.385
$ws.Range("N32").Value = -22793.928
# Row 45
$ws.Range("H45").Value = 2734.7097
$ws.Range("I45").Value = 2490.3809
$ws.Range("J45").Value = 3247.8
$ws.Range("K45").Value = 2490.3809
$ws.Range("L45").Value = 3247.8
$ws.Range("M45").Value = -2113.3809
$ws.Range("N45").Value = -4001.8
# Row 102
$ws.Range("H102").Value = 1524.5454
$ws.Range("I102").Value = 1524.5454
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1524.5454
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 97.45460000000003
$ws.Range("N102").ClearContents()
# Row 116
$ws.Range("H116").Value = 1130.2084
$ws.Range("I116").Value = 1095.8235
$ws.Range("J116").Value = 1213.7142
$ws.Range("K116").Value = 1095.8235
$ws.Range("L116").Value = 1213.7142
$ws.Range("M116").Value = 1198.1765
$ws.Range("N116").Value = -5801.7142
# Row 122
$ws.Range("H122").Value = 2822.7896
$ws.Range("I122").Value = 2765.875
$ws.Range("K122").Value = 8297.625
$ws.Range("M122").Value = -5847.625
# Row 125
$ws.Range("H125").Value = 34898
$ws.Range("J125").Value = 34898
$ws.Range("L125").Value = 34898
$ws.Range("N125").Value = -44738
# Row 139
$ws.Range("H139").Value = 48810
$ws.Range("J139").Value = 48810
$ws.Range("L139").Value = 48810
$ws.Range("N139").Value = -59090

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1130.2084
$ws.Range("I3").Value = 1095.8235
$ws.Range("J3").Value = 1213.7142
$ws.Range("K3").Value = 1095.8235
$ws.Range("L3").Value = 1213.7142
$ws.Range("M3").Value = -981.8235
$ws.Range("N3").Value = -1441.7142
# Row 99
$ws.Range("H99").Value = 1280
$ws.Range("I99").Value = 1300
$ws.Range("K99").Value = 1300
$ws.Range("M99").Value = 198
# Row 135
$ws.Range("H135").Value = 39492.832
$ws.Range("J135").Value = 39492.832
$ws.Range("L135").Value = 39492.832
$ws.Range("N135").Value = -49632.832

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
# Row 141
$ws.Range("H141").Value = 24649.455
$ws.Range("J141").Value = 24649.455
$ws.Range("L141").Value = 24649.455
$ws.Range("N141").Value = -35009.455

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 52.22222
$ws.Range("I33").Value = 18.333334
$ws.Range("K33").Value = 110.000004
$ws.Range("M33").Value = 172.999996
# Row 63
$ws.Range("H63").Value = 2779.3333
$ws.Range("I63").Value = 1500
$ws.Range("K63").Value = 4500
$ws.Range("M63").Value = -3751
# Row 66
$ws.Range("H66").Value = 2779.3333
$ws.Range("I66").Value = 1500
$ws.Range("K66").Value = 13500
$ws.Range("M66").Value = -9756
# Row 75
$ws.Range("H75").Value = 450
$ws.Range("J75").Value = 450
$ws.Range("L75").Value = 1350
$ws.Range("N75").Value = -3346
# Row 78
$ws.Range("H78").Value = 450
$ws.Range("J78").Value = 450
$ws.Range("L78").Value = 4050
$ws.Range("N78").Value = -14034
# Row 80
$ws.Range("H80").Value = 22400.2
$ws.Range("J80").Value = 36000.332
$ws.Range("L80").Value = 108000.996
$ws.Range("N80").Value = -109872.996
# Row 83
$ws.Range("H83").Value = 22400.2
$ws.Range("J83").Value = 36000.332
$ws.Range("L83").Value = 324002.988
$ws.Range("N83").Value = -333362.988
# Row 107
$ws.Range("H107").Value = 6051.697
$ws.Range("I107").Value = 8544
$ws.Range("J107").Value = 319.4
$ws.Range("K107").Value = 25632
$ws.Range("L107").Value = 958.1999999999999
$ws.Range("M107").Value = -23712
$ws.Range("N107").Value = -4798.2
# Row 112
$ws.Range("H112").Value = 100002080
$ws.Range("I112").Value = 1141.6666
$ws.Range("J112").Value = 250003490
$ws.Range("K112").Value = 3424.9998
$ws.Range("L112").Value = 750010470
$ws.Range("M112").Value = -2316.9998
$ws.Range("N112").Value = -750012686
# Row 117
$ws.Range("H117").Value = 1699.091
$ws.Range("J117").Value = 2933
$ws.Range("L117").Value = 8799
$ws.Range("N117").Value = -15683
# Row 131
$ws.Range("H131").Value = 756.0700000000001
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 756.0700000000001
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2268.21
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12348.21
# Row 137
$ws.Range("H137").Value = 18524412
$ws.Range("J137").Value = 22228994
$ws.Range("L137").Value = 66686982
$ws.Range("N137").Value = -66697182

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = 0
# Row 102
$ws.Range("H102").Value = 6335.4287
$ws.Range("I102").Value = 5722.3335
$ws.Range("J102").Value = 10014
$ws.Range("K102").Value = 5722.3335
$ws.Range("L102").Value = 10014
$ws.Range("M102").Value = -4100.3335
$ws.Range("N102").Value = -13258
# Row 107
$ws.Range("H107").Value = 561.25
$ws.Range("I107").Value = 453.66666
$ws.Range("J107").Value = 668.8333
$ws.Range("K107").Value = 453.66666
$ws.Range("L107").Value = 668.8333
$ws.Range("M107").Value = 1466.33334
$ws.Range("N107").Value = -4508.8333
# Row 109
$ws.Range("H109").Value = 28257
$ws.Range("J109").Value = 28257
$ws.Range("L109").Value = 28257
$ws.Range("N109").Value = -30337

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1099.5
$ws.Range("I46").Value = 999
$ws.Range("K46").Value = 999
$ws.Range("M46").Value = -811

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 18000
$ws.Range("J54").Value = 18000
$ws.Range("L54").Value = 18000
$ws.Range("N54").Value = -19040
# Row 75
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
# Row 78
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360
# Row 101
$ws.Range("H101").Value = 17900
$ws.Range("J101").Value = 17900
$ws.Range("L101").Value = 17900
$ws.Range("N101").Value = -24390
# Row 113
$ws.Range("H113").Value = 1952.4445
$ws.Range("I113").Value = 2357.7693
$ws.Range("J113").Value = 898.6
$ws.Range("K113").Value = 7073.3079
$ws.Range("L113").Value = 2695.8
$ws.Range("M113").Value = -4903.3079
$ws.Range("N113").Value = -7035.8
# Row 126
$ws.Range("H126").Value = 2214.3157
$ws.Range("J126").Value = 3501
$ws.Range("L126").Value = 10503
$ws.Range("N126").Value = -15443
